$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.788.98"
$ws.Range("E2").Value = "  +12.01%  "

$ws.Range("D3").Value = "1.685.41"
$ws.Range("E3").Value = "  +6.38%  "

$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "'306.33"
$ws.Range("E5").Value = "  +3.02%  "

$ws.Range("D6").Value = "'0.9960"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").Value = "'0.3678"
$ws.Range("E7").Value = "  +2.20%  "

$ws.Range("D8").Value = "'50.25"
$ws.Range("E8").Value = "  +22.65%  "

$ws.Range("D9").Value = "'0.3419"
$ws.Range("E9").Value = "  +3.09%  "

$ws.Range("D10").Value = "'1.160"
$ws.Range("E10").Value = "  +4.89%  "

$ws.Range("D11").Value = "'0.07214"
$ws.Range("E11").Value = "  +4.63%  "

$ws.Range("D12").Value = "'0.9973"
$ws.Range("E12").Value = "  -0.40%  "

$ws.Range("D13").Value = "'6.099"
$ws.Range("E13").Value = "  +5.58%  "

$ws.Range("D14").Value = "'20.13"
$ws.Range("E14").Value = "  +4.92%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.682.38"
$ws.Range("E15").Value = "  +6.15%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'6.672"
$ws.Range("E16").Value = "  +3.04%  "

$ws.Range("D17").Value = "'0.00001102"
$ws.Range("E17").Value = "  +4.25%  "

$ws.Range("D18").Value = "'0.9957"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("D19").Value = "'0.06658"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").Value = "'80.93"
$ws.Range("E20").Value = "  +7.21%  "

$ws.Range("D21").Value = "'16.31"
$ws.Range("E21").Value = "  +3.83%  "

$ws.Range("D22").Value = "'6.059"
$ws.Range("E22").Value = "  +3.33%  "

$ws.Range("E23").Value = "  +5.34%  "

$ws.Range("D24").Value = "24.694.46"
$ws.Range("E24").Value = "  +11.56%  "

$ws.Range("D25").Value = "'2.409"
$ws.Range("E25").Value = "  +2.10%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.651"
$ws.Range("E26").Value = "  +6.73%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'152.54"
$ws.Range("E27").Value = "  +3.34%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.41"
$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "1.868.74"
$ws.Range("E29").Value = "  +6.50%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'127.11"
$ws.Range("E30").Value = "  +5.19%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'6.221"
$ws.Range("E31").Value = "  +7.10%  "

$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "'4.024"
$ws.Range("E32").Value = "  +2.86%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.9814"
$ws.Range("E33").Value = "  +7.80%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.08381"
$ws.Range("E34").Value = "  +3.45%  "

$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.689"
$ws.Range("E35").Value = "  +4.05%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'12.30"
$ws.Range("E36").Value = "  +6.38%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06341"
$ws.Range("E37").Value = "  +6.45%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.285"
$ws.Range("E38").Value = "  +4.74%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02299"
$ws.Range("E39").Value = "  +6.19%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.606"
$ws.Range("E40").Value = "  +2.88%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.235"
$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.2082"
$ws.Range("E42").Value = "  +5.93%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.6081"
$ws.Range("E43").Value = "  +6.27%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'0.9955"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.759"
$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("D46").Value = "'13.03"
$ws.Range("E46").Value = "  +3.76%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5866"
$ws.Range("E47").Value = "  +6.12%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'125.14"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.992"
$ws.Range("E49").Value = "  +3.59%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07234"
$ws.Range("E50").Value = "  +8.25%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'75.49"
$ws.Range("E51").Value = "  +5.16%  "
